$d = $word.ActiveDocument

# --- Change 1: remove the existing "_GoBack" bookmark after "...reg@fjk.hu címre." ---
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

# --- Change 2: split the run containing "...megrendelt szolgáltatások, a szállást..."
#     right after "szolgáltatáso" and insert a new "_GoBack" bookmark at the split point ---
$r2 = $d.Content
$found2 = $r2.Find.Execute("szolgáltatáso", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found2) {
    $r2.Collapse(0)
    $d.Bookmarks.Add("_GoBack", $r2)
}

# --- Change 3: rewrite the bold red paragraph about foreign bank transfers ---
$r3 = $d.Content
$oldText3 = "Amennyiben az utalást csak külföldi (nem magyar) bankszámláról tudnád intézni, úgy nem vagy köteles az előleg fizetésére tekintettel a magas határon kívüli utalási költségekre. A szállás költségét a Benczúr Hotelben tudod majd egyenlíteni, a regisztráció fennmaradó költséget pedig a konferencia regisztrációs asztalánál készpénzben. Ebben az esetben egy emailt fogunk küldeni a konferencia előtt, melyben megkérünk majd egy második visszaigazolásra a megrendelt szolgáltatásokról."
$found3 = $r3.Find.Execute($oldText3, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found3) {
    $r3.Text = "Amennyiben az utalást csak külföldi (nem magyar) bankszámláról tudnád intézni, úgy nem vagy köteles az előleg fizetésére tekintettel a magas határon kívüli utalási költségekre. A regisztráció költségét a konferencián személyesen, készpénzben tudod egyenlíteni. Ebben az esetben egy e-mailt fogunk küldeni a konferencia előtt, melyben megkérünk majd egy második visszaigazolásra a megrendelt szolgáltatásokról."
}
